$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A (Id), Q (Ost), R (Nord) per row.
# These represent a cyclic re-shuffle of the original row values.
$updates = @{
    2  = @(111336408, 625192.9989858982, 7209753.182408583)
    3  = @(111336410, 625259.9138955096, 7209755.994336623)
    4  = @(111336414, 625154.518985283,  7209790.550712772)
    5  = @(111336398, 625080.5973707421, 7209863.719748351)
    6  = @(111336416, 625105.0111505401, 7209827.879692691)
    7  = @(111336412, 625184.8249035137, 7209765.975211025)
    8  = @(111336419, 625180.6727454782, 7209784.021863313)
    10 = @(111336404, 625133.9709027896, 7209754.516890368)
    11 = @(111336415, 625148.1990682605, 7209799.60762905)
    12 = @(111336409, 625220.7730415409, 7209758.586850428)
    13 = @(111336413, 625162.5685057295, 7209780.71893465)
    14 = @(111336405, 625163.9543035047, 7209747.724855823)
    16 = @(111336403, 625077.0981605061, 7209815.690764531)
    17 = @(111336407, 625184.0874587877, 7209753.231787121)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 1).Value2 = $vals[0]
    $ws.Cells.Item($row, 17).Value2 = $vals[1]
    $ws.Cells.Item($row, 18).Value2 = $vals[2]
}
